# Apply "authors' replies" updates to the provenance/versioning comparison tables.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Header row 2: rename the "File" sub-header under "Distribution" (column G) to "Local"
    # so it is distinguishable from the duplicate "File" sub-header under "Storage" (column F).
    if ($ws.Range("G2").Value2 -eq "File") {
        $ws.Range("G2").Value2 = "Local"
    }
}

# "script" sheet (sheet1) specific data corrections based on authors' replies.
$ws1 = $wb.Worksheets.Item("script")

# Datatrack: Versioning -> "Trial ID" (was "Sequence")
$ws1.Range("I8").Value2 = "Trial ID"

# RDataTracker: Artifacts -> "PROV-JSON" (was "Interoperable (PROV)")
$ws1.Range("C18").Value2 = "PROV-JSON"
# RDataTracker: Versioning -> "Trial ID" (was "Sequence")
$ws1.Range("I18").Value2 = "Trial ID"

# SPADE: Artifacts -> expanded list of supported backends
$ws1.Range("C21").Value2 = "PostgreSQL, MySQL, H2, Neo4j, Datalog, GraphViz, PROV"

# versuchung: Artifacts -> add SQLite to the list
$ws1.Range("C26").Value2 = "Content DB, SQLite, Proprietary (Dict)"
